# Insert a new weekly price observation as row 196 in the "Hortaliza,
# Femacal de La Calera - Albahaca" sheet, pushing the existing rows
# 196-208 down to 197-209 (dimension grows from A1:R208 to A1:R209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 196..208 down by one, carrying formats along (mirrors
# Excel's own Insert behaviour of copying the row above's formatting).
$ws.Rows("196:196").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A196").Value = 3
$ws.Range("B196").Value = "Femacal de La Calera"
$ws.Range("C196").Value = "Coquimbo"
$ws.Range("D196").Value = 44931
$ws.Range("E196").Value = 5
$ws.Range("F196").Value = 100112052
$ws.Range("G196").Value = "Albahaca"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 130
$ws.Range("K196").Value = 5000
$ws.Range("L196").Value = 5500
$ws.Range("M196").Value = 5192
$ws.Range("N196").Value = "`$/docena de matas"
$ws.Range("O196").Value = "Provincia de Quillota"
$ws.Range("P196").Value = 865
$ws.Range("Q196").Value = 6
$ws.Range("R196").Value = "Hortaliza"
